# Refresh crypto price/volume data (and two pairs of swapped rows) to match the
# latest GitHub Actions scrape, per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '64.295.15'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  -2.83%  '
# Row 3
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.169.94'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  -4.40%  '
# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '569.64'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  -2.72%  '
# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '168.13'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  -8.30%  '
# Row 7
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.602'
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  -7.27%  '
# Row 8
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '1.00'
# Row 9
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '3.177.03'
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  -4.15%  '
# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.120'
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  -4.20%  '
# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '6.79'
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  -0.52%  '
# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.385'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  -3.86%  '
# Row 13
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '3.730.61'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  -4.18%  '
# Row 14
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  -2.08%  '
# Row 15
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '64.381.77'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  -2.77%  '
# Row 16
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '25.30'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  -3.05%  '
# Row 17
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.0000158'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  -2.99%  '
# Row 18
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '3.176.64'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  -3.45%  '
# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '416.69'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  -1.77%  '
# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '12.90'
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  -2.16%  '
# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '5.34'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  -3.59%  '
# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '7.08'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  -4.05%  '
# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  -0.08%  '
# Row 24
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '70.31'
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  -2.04%  '
# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '5.67'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  -0.16%  '
# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.204'
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  +1.65%  '
# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.486'
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  -5.20%  '
# Row 28
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.0000105'
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  -7.16%  '
# Row 29
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '8.73'
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  -2.00%  '
# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.999'
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  -0.12%  '
# Row 31
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  -3.52%  '
# Row 32
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '21.68'
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  -3.12%  '
# Row 33
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  -0.12%  '
# Row 34
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '5.04'
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  -2.44%  '
# Row 35
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '6.31'
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  -3.87%  '
# Row 36
$ws.Range('B36').Value = 'Monero'
$ws.Range('C36').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '158.03'
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  -1.42%  '
# Row 37
$ws.Range('B37').Value = 'Fetch.AI'
$ws.Range('C37').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.13'
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  -4.38%  '
# Row 38
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.35'
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  -5.74%  '
# Row 39
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.708.70'
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  -6.02%  '
# Row 40
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.70'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  -5.52%  '
# Row 41
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '24.28'
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  -7.95%  '
# Row 42
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '4.19'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  -2.74%  '
# Row 43
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '39.13'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  -2.06%  '
# Row 44
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.714'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  -6.46%  '
# Row 45
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0621'
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  -6.43%  '
# Row 46
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '5.58'
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  -5.92%  '
# Row 47
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0263'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  -3.12%  '
# Row 48
$ws.Range('B48').Value = 'Bittensor'
$ws.Range('C48').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '292.35'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  -6.85%  '
# Row 49
$ws.Range('B49').Value = 'InjectiveProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '21.39'
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  -7.57%  '
# Row 50
$ws.Range('B50').Value = 'FirstDigitalUSD'
$ws.Range('C50').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.998'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  -0.17%  '
# Row 51
$ws.Range('B51').Value = 'dogwifhat'
$ws.Range('C51').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.00'
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  -12.70%  '
